# Rearranges the weekly price rows (rows 2-12) on the active sheet.
# The columns that vary per-row are: D (Fecha), M (Volumen), N (Precio minimo),
# O (Precio maximo), P (Precio promedio ponderado), S (Precio $/Kg).
# Each target row ends up with the values that used to live at its "source" row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# target row -> source row (i.e. target row receives source row's old values)
$rowMap = @{
    2  = 6
    3  = 2
    4  = 5
    5  = 9
    6  = 3
    7  = 11
    8  = 12
    9  = 8
    10 = 10
    11 = 4
    12 = 7
}

$cols = @(4, 13, 14, 15, 16, 19)  # D, M, N, O, P, S

# Snapshot the original values for every row/column we touch before writing
# anything, since several target rows also act as sources for other targets.
$snapshot = @{}
foreach ($col in $cols) {
    for ($row = 2; $row -le 12; $row++) {
        $snapshot["$row-$col"] = $ws.Cells.Item($row, $col).Value2
    }
}

foreach ($targetRow in $rowMap.Keys) {
    $sourceRow = $rowMap[$targetRow]
    foreach ($col in $cols) {
        $ws.Cells.Item($targetRow, $col).Value = $snapshot["$sourceRow-$col"]
    }
}
